$wb = $excel.ActiveWorkbook

# --- Rename sheets (Monthly -> Weekly) ---
$wsGeneral = $wb.Worksheets.Item("GeneralTaxRateMonthly")
$wsGeneral.Name = "GeneralTaxRateWeekly"

$wsProcess = $wb.Worksheets.Item("ProcessPayrollForMonthlyTax")
$wsProcess.Name = "ProcessPayrollForWeeklyTax"

$wsFirst = $wb.Worksheets.Item("first")
$wsReports = $wb.Worksheets.Item("TestReports")

# --- Update cell content to reflect the Weekly module ---
$wsFirst.Range("A3").Value = "GeneralTaxRateWeekly"
$wsFirst.Range("A4").Value = "ProcessPayrollForWeeklyTax"

$wsGeneral.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 107"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"
$wsReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"

# --- Update each sheet's own stored selection/scroll position ---
[void]$wsGeneral.Activate()
[void]$wsGeneral.Range("A2").Select()

[void]$wsProcess.Activate()
[void]$wsProcess.Range("B9").Select()

[void]$wsReports.Activate()
[void]$wsReports.Range("B2").Select()

# --- Make "first" the active sheet/tab, matching the target workbook state ---
[void]$wsFirst.Activate()
[void]$wsFirst.Range("A3").Select()
